# Settings back to qfin22 presentation - recalibration
#
# Updates the calibration results on the "Linear" and "NonLinear" sheets,
# and the matching "abs_epsi_autocorr" series stored alongside each, to the
# freshly recalibrated values.

$wb = $excel.ActiveWorkbook

# --- Linear sheet (param / mu / B / sig2 / abs_epsi_autocorr) ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = 0.00003073738816579158
$wsLinear.Range("B3").Value = -0.17783286840234549531
$wsLinear.Range("B4").Value = 0.00051293755516057494
$wsLinear.Range("B5").Value = "[0.9999999999999999, 0.33852263214045586, 0.2404927462155008, 0.359046808988458, 0.23394170318712887, 0.24314402051240827, 0.390243582803775, 0.44997256760421533, 0.30320666092308496, 0.23828761067604057, 0.21258249725642994, 0.2919246191822665, 0.25781010706010127, 0.3358444684246943, 0.45853550221705297, 0.2904081765684354, 0.1727121050159306, 0.27171992978770526, 0.21050918218148212, 0.21889085079828266]"

# --- NonLinear sheet (c / p / mu_0 / B_0 / sig2_0 / mu_1 / B_1 / sig2_1 / abs_epsi_autocorr) ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = 0.00029706998665980189
$wsNonLinear.Range("B5").Value = -0.11131810811062360644
$wsNonLinear.Range("B6").Value = 0.00033623157484431818
$wsNonLinear.Range("B7").Value = 0.00013557436067022311
$wsNonLinear.Range("B8").Value = -0.22962644471145790548
$wsNonLinear.Range("B9").Value = 0.00067870522623198414
$wsNonLinear.Range("B10").Value = "[1.0, 0.34455119405311185, 0.24523926409817196, 0.3542710037954998, 0.23255401732804715, 0.24330078992097753, 0.3903375579713773, 0.4491816443403333, 0.3006175022251547, 0.23866163314023192, 0.2127588019306442, 0.2892005604967266, 0.25584119682246625, 0.3374641466445282, 0.45582917456726163, 0.29005969926464864, 0.1719881393727604, 0.26795743826021057, 0.2104805233905603, 0.21854800642370004]"
